$wb = $excel.ActiveWorkbook

# Sheet references
$ws2 = $wb.Worksheets.Item("630-845 AM")   # xl/worksheets/sheet2.xml
$ws3 = $wb.Worksheets.Item("900-1145 AM")  # xl/worksheets/sheet3.xml

# Populate the "630-845 AM" tab with sample counts (B2:H11)
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 19
$ws2.Range("D2").Value = 200
$ws2.Range("E2").Value = 30
$ws2.Range("F2").Value = 1
$ws2.Range("H2").Value = 11
$ws2.Range("B3").Value = 11
$ws2.Range("C3").Value = 18
$ws2.Range("E3").Value = 40
$ws2.Range("F3").Value = 2
$ws2.Range("G3").Value = 1
$ws2.Range("H3").Value = 11
$ws2.Range("B4").Value = 22
$ws2.Range("C4").Value = 17
$ws2.Range("D4").Value = 300
$ws2.Range("F4").Value = 3
$ws2.Range("G4").Value = 1
$ws2.Range("H4").Value = 11
$ws2.Range("B5").Value = 33
$ws2.Range("C5").Value = 16
$ws2.Range("D5").Value = 400
$ws2.Range("E5").Value = 60
$ws2.Range("F5").Value = 4
$ws2.Range("G5").Value = 1
$ws2.Range("H5").Value = 11
$ws2.Range("B6").Value = 44
$ws2.Range("C6").Value = 15
$ws2.Range("D6").Value = 500
$ws2.Range("E6").Value = 70
$ws2.Range("F6").Value = 5
$ws2.Range("G6").Value = 1
$ws2.Range("H6").Value = 11
$ws2.Range("B7").Value = " "
$ws2.Range("C7").Value = 14
$ws2.Range("D7").Value = 600
$ws2.Range("E7").Value = 80
$ws2.Range("G7").Value = 1
$ws2.Range("H7").Value = 11
$ws2.Range("B8").Value = 66
$ws2.Range("C8").Value = 12
$ws2.Range("D8").Value = 700
$ws2.Range("E8").Value = 90
$ws2.Range("F8").Value = 7
$ws2.Range("G8").Value = 1
$ws2.Range("H8").Value = " "
$ws2.Range("B9").Value = 77
$ws2.Range("C9").Value = 11
$ws2.Range("D9").Value = 800
$ws2.Range("E9").Value = 100
$ws2.Range("F9").Value = 8
$ws2.Range("G9").Value = 1
$ws2.Range("H9").Value = 11
$ws2.Range("B10").Value = 88
$ws2.Range("D10").Value = 900
$ws2.Range("E10").Value = 110
$ws2.Range("F10").Value = 9
$ws2.Range("G10").Value = 1
$ws2.Range("H10").Value = 11
$ws2.Range("B11").Value = 99
$ws2.Range("C11").Value = 9
$ws2.Range("D11").Value = 1000
$ws2.Range("E11").Value = 120
$ws2.Range("F11").Value = 9
$ws2.Range("G11").Value = 1
$ws2.Range("H11").Value = 11

# Move the cursor on the "900-1145 AM" tab (A12 -> A14), without making it the active tab
$ws3.Range("A14").Select() | Out-Null

# Make "630-845 AM" the active sheet/tab and move its cursor to H12
$ws2.Activate()
$ws2.Range("H12").Select() | Out-Null
